$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "283.10"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "2.08%"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "28.57"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "4.72%"
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "5.066"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "3.73%"
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.06483"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "1.03%"
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "7.207"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "3.48%"
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.371"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "16.46%"
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.9110"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "2.93%"
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "0.02%"
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.06527"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "27.06%"
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07608"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "2.23%"
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.02753"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "-5.06%"
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.08947"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "-0.35%"
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "1.57%"
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.0006356"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "-0.31%"
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.006072"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "-1.54%"
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "-0.73%"
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.367"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "1.37%"
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "-1.45%"
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "1.17%"
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.988"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "2.03%"
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.1564"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "4.18%"
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.04454"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "0.51%"
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.001181"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "0.36%"
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.004355"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "12.67%"
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "1.74%"
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.0001637"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "-15.68%"
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.04112"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "-0.87%"
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.006645"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "-2.04%"
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1230"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "4.79%"
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.002073"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "7.30%"
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.01246"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "8.56%"
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00005400"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "1.61%"
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "-0.16%"
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.933"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "14.55%"
$ws.Cells.Item(47, 5).Style = "Normal"
